$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price reports was added at the top of the data table
# (rows 33-34), pushing every subsequent row down by two positions.
$ws.Rows("33:34").Insert()

# New row 33: Sandia, Primera, week of 2022-02-09 (serial 44601)
$ws.Range("A33").Value = 11
$ws.Range("B33").Value = "Vega Monumental Concepción"
$ws.Range("C33").Value = "Bíobío"
$ws.Range("D33").Value = 44601
$ws.Range("E33").Value = 8
$ws.Range("F33").Value = 100112028
$ws.Range("G33").Value = "Sandia"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 800
$ws.Range("K33").Value = 2000
$ws.Range("L33").Value = 2500
$ws.Range("M33").Value = 2188
$ws.Range("N33").Value = "$/unidad"
$ws.Range("O33").Value = "Región Metropolitana"
$ws.Range("P33").Value = 2188
$ws.Range("Q33").Value = 1
$ws.Range("R33").Value = "Hortaliza"

# New row 34: Sandia, Segunda, same week of 2022-02-09 (serial 44601)
$ws.Range("A34").Value = 11
$ws.Range("B34").Value = "Vega Monumental Concepción"
$ws.Range("C34").Value = "Bíobío"
$ws.Range("D34").Value = 44601
$ws.Range("E34").Value = 8
$ws.Range("F34").Value = 100112028
$ws.Range("G34").Value = "Sandia"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Segunda"
$ws.Range("J34").Value = 700
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = 1786
$ws.Range("N34").Value = "$/unidad"
$ws.Range("O34").Value = "Región Metropolitana"
$ws.Range("P34").Value = 1786
$ws.Range("Q34").Value = 1
$ws.Range("R34").Value = "Hortaliza"
